$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3781.8474
$ws.Range("I137").Value = 4951.0884
$ws.Range("J137").Value = 2191.68
$ws.Range("K137").Value = 14853.2652
$ws.Range("L137").Value = 6575.039999999999
$ws.Range("M137").Value = -12303.2652
$ws.Range("N137").Value = -11675.04
$ws.Range("H138").Value = 2291.1558
$ws.Range("J138").Value = 3551.647
$ws.Range("L138").Value = 10654.941
$ws.Range("N138").Value = -20934.941

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4102.9434
$ws.Range("I32").Value = 3662.9048
$ws.Range("J32").Value = 13343.75
$ws.Range("K32").Value = 3662.9048
$ws.Range("L32").Value = 13343.75
$ws.Range("M32").Value = -3375.9048
$ws.Range("N32").Value = -13917.75
$ws.Range("H61").Value = 3871.9648
$ws.Range("I61").Value = 4107.2246
$ws.Range("J61").Value = 2431
$ws.Range("K61").Value = 4107.2246
$ws.Range("L61").Value = 2431
$ws.Range("M61").Value = -3895.2246
$ws.Range("N61").Value = -2855
$ws.Range("H74").Value = 1990.92
$ws.Range("I74").Value = 1120.5652
$ws.Range("K74").Value = 1120.5652
$ws.Range("M74").Value = -246.5652
$ws.Range("H77").Value = 1990.92
$ws.Range("I77").Value = 1120.5652
$ws.Range("K77").Value = 5602.826
$ws.Range("M77").Value = -1234.826
$ws.Range("H96").Value = 29538.23
$ws.Range("J96").Value = 29538.23
$ws.Range("L96").Value = 29538.23
$ws.Range("N96").Value = -35030.23
$ws.Range("H101").Value = 36998.25
$ws.Range("J101").Value = 36998.25
$ws.Range("L101").Value = 36998.25
$ws.Range("N101").Value = -43488.25
$ws.Range("H136").Value = 3871.9648
$ws.Range("I136").Value = 4107.2246
$ws.Range("J136").Value = 2431
$ws.Range("K136").Value = 12321.6738
$ws.Range("L136").Value = 7293
$ws.Range("M136").Value = -9771.673799999999
$ws.Range("N136").Value = -12393

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3751.5264
$ws.Range("I20").Value = 3753.9333
$ws.Range("K20").Value = 3753.9333
$ws.Range("M20").Value = -3506.9333
$ws.Range("H134").Value = 4031.5557
$ws.Range("I134").Value = 3169.261
$ws.Range("K134").Value = 9507.782999999999
$ws.Range("M134").Value = -6972.782999999999

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1855.96
$ws.Range("I31").Value = 1557.8445
$ws.Range("J31").Value = 2303.1333
$ws.Range("K31").Value = 1557.8445
$ws.Range("L31").Value = 2303.1333
$ws.Range("M31").Value = -1262.8445
$ws.Range("N31").Value = -2893.1333
$ws.Range("H34").Value = 1855.96
$ws.Range("I34").Value = 1557.8445
$ws.Range("J34").Value = 2303.1333
$ws.Range("K34").Value = 1557.8445
$ws.Range("L34").Value = 2303.1333
$ws.Range("M34").Value = -1355.8445
$ws.Range("N34").Value = -2707.1333
$ws.Range("H58").Value = 1358.8684
$ws.Range("I58").Value = 1314.5135
$ws.Range("K58").Value = 1314.5135
$ws.Range("M58").Value = -1111.5135
$ws.Range("H132").Value = 5562.574
$ws.Range("I132").Value = 2003.5581
$ws.Range("K132").Value = 6010.6743
$ws.Range("M132").Value = -3480.6743
$ws.Range("H134").Value = 5440.85
$ws.Range("I134").Value = 6109.25
$ws.Range("J134").Value = 2767.25
$ws.Range("K134").Value = 18327.75
$ws.Range("L134").Value = 8301.75
$ws.Range("M134").Value = -15792.75
$ws.Range("N134").Value = -13371.75
$ws.Range("H136").Value = 1358.8684
$ws.Range("I136").Value = 1314.5135
$ws.Range("K136").Value = 3943.5405
$ws.Range("M136").Value = -1393.5405

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1697.409
$ws.Range("I5").Value = 1250.2106
$ws.Range("J5").Value = 4529.6665
$ws.Range("K5").Value = 3750.6318
$ws.Range("L5").Value = 13588.9995
$ws.Range("M5").Value = -3638.6318
$ws.Range("N5").Value = -13812.9995
$ws.Range("H31").Value = 4500
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 4500
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 13500
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -14076
$ws.Range("H97").Value = 650.63635
$ws.Range("I97").Value = 326.33334
$ws.Range("K97").Value = 979.0000200000001
$ws.Range("M97").Value = -483.0000200000001
$ws.Range("H135").Value = 1697.409
$ws.Range("I135").Value = 1250.2106
$ws.Range("J135").Value = 4529.6665
$ws.Range("K135").Value = 11251.8954
$ws.Range("L135").Value = 40766.9985
$ws.Range("M135").Value = -8716.895400000001
$ws.Range("N135").Value = -45836.9985

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 238393.47
$ws.Range("I3").Value = 1501500
$ws.Range("J3").Value = 8737.727999999999
$ws.Range("K3").Value = 1501500
$ws.Range("L3").Value = 8737.727999999999
$ws.Range("M3").Value = -1501384
$ws.Range("N3").Value = -8969.727999999999
$ws.Range("H11").Value = 1264962.9
$ws.Range("I11").Value = 3339999.8
$ws.Range("K11").Value = 3339999.8
$ws.Range("M11").Value = -3339860.8
$ws.Range("H80").Value = 36925730
$ws.Range("I80").Value = 53334810
$ws.Range("J80").Value = 5299.5
$ws.Range("K80").Value = 53334810
$ws.Range("L80").Value = 5299.5
$ws.Range("M80").Value = -53333812
$ws.Range("N80").Value = -7295.5
$ws.Range("H83").Value = 36925730
$ws.Range("I83").Value = 53334810
$ws.Range("J83").Value = 5299.5
$ws.Range("K83").Value = 266674050
$ws.Range("L83").Value = 26497.5
$ws.Range("M83").Value = -266669058
$ws.Range("N83").Value = -36481.5
$ws.Range("H132").Value = 9044.451999999999
$ws.Range("I132").Value = 7117.5806
$ws.Range("J132").Value = 14474.728
$ws.Range("K132").Value = 21352.7418
$ws.Range("L132").Value = 43424.18399999999
$ws.Range("M132").Value = -18822.7418
$ws.Range("N132").Value = -48484.18399999999

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 12500662
$ws.Range("I16").Value = 12500662
$ws.Range("K16").Value = 12500662
$ws.Range("M16").Value = -12500492
$ws.Range("H33").Value = 17500
$ws.Range("I33").Value = 17500
$ws.Range("K33").Value = 17500
$ws.Range("M33").Value = -17210
$ws.Range("H132").Value = 83111.7
$ws.Range("I132").Value = 87222.84
$ws.Range("K132").Value = 261668.52
$ws.Range("M132").Value = -259138.52
$ws.Range("H136").Value = 3755054.2
$ws.Range("I136").Value = 5629895
$ws.Range("K136").Value = 16889685
$ws.Range("M136").Value = -16887135

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 2250
$ws.Range("I18").Value = 2250
$ws.Range("K18").Value = 2250
$ws.Range("M18").Value = -2077
$ws.Range("H26").Value = 63399.57
$ws.Range("I26").Value = 71933.336
$ws.Range("K26").Value = 71933.336
$ws.Range("M26").Value = -71640.336
$ws.Range("H38").Value = 11068.5
$ws.Range("I38").Value = 9482.200000000001
$ws.Range("K38").Value = 9482.200000000001
$ws.Range("M38").Value = -9009.200000000001
$ws.Range("H81").Value = 5053571.5
$ws.Range("I81").Value = 6495864
$ws.Range("J81").Value = 5547.5
$ws.Range("K81").Value = 12991728
$ws.Range("L81").Value = 11095
$ws.Range("M81").Value = -12990667
$ws.Range("N81").Value = -13217
$ws.Range("H84").Value = 5053571.5
$ws.Range("I84").Value = 6495864
$ws.Range("J84").Value = 5547.5
$ws.Range("K84").Value = 64958640
$ws.Range("L84").Value = 55475
$ws.Range("M84").Value = -64953336
$ws.Range("N84").Value = -66083
$ws.Range("H122").Value = 4002.0322
$ws.Range("I122").Value = 3388.8635
$ws.Range("J122").Value = 5500.8887
$ws.Range("K122").Value = 10166.5905
$ws.Range("L122").Value = 16502.6661
$ws.Range("M122").Value = -7716.5905
$ws.Range("N122").Value = -21402.6661
$ws.Range("H132").Value = 3793.2122
$ws.Range("I132").Value = 3480.4546
$ws.Range("J132").Value = 4418.727
$ws.Range("K132").Value = 10441.3638
$ws.Range("L132").Value = 13256.181
$ws.Range("M132").Value = -7911.363799999999
$ws.Range("N132").Value = -18316.181
$ws.Range("H136").Value = 1701.4127
$ws.Range("I136").Value = 1783.3019
$ws.Range("K136").Value = 5349.905699999999
$ws.Range("M136").Value = -2799.905699999999
